$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 8-10 (sending cluster = ECs), new data has only FAPs/MuSCs as senders
$ws.Range("A8:T10").Delete()

# Update rows 2-7 with the new TPM-derived values (FAPs/MuSCs sending Fgf18 -> Fgfr1)
$row2 = @("FAPs", "Fgf18", "Fgfr1", "ECs", 3.0, 1.0, 10.71557066666667, 32.146712, 0.9375025736567436, 0.9375025736567436, 3.0, 1.0, 1.845768666666667, 5.537306, 0.01459089321241885, 0.01459089321241885, 19.77846458198578, 178.006181237872, 0.01367899993859338, 0.01367899993859338)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(2, $c).Value = $row2[$c - 1] }

$row3 = @("FAPs", "Fgf18", "Fgfr1", "FAPs", 3.0, 1.0, 10.71557066666667, 32.146712, 0.9375025736567436, 0.9375025736567436, 3.0, 1.0, 82.95722966666666, 248.871689, 0.6557810310272387, 0.6557810310272387, 888.9340568040631, 8000.406511236568, 0.6147964043433092, 0.6147964043433092)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(3, $c).Value = $row3[$c - 1] }

$row4 = @("FAPs", "Fgf18", "Fgfr1", "MuSCs", 3.0, 1.0, 10.71557066666667, 32.146712, 0.9375025736567436, 0.9375025736567436, 3.0, 1.0, 41.69841866666667, 125.095256, 0.3296280757603424, 0.3296280757603424, 446.8223519109191, 4021.401167198272, 0.3090271693748411, 0.3090271693748411)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(4, $c).Value = $row4[$c - 1] }

$row5 = @("MuSCs", "Fgf18", "Fgfr1", "ECs", 2.0, 0.6666666666666666, 0.71434, 2.14302, 0.06249742634325634, 0.06249742634325634, 3.0, 1.0, 1.845768666666667, 5.537306, 0.01459089321241885, 0.01459089321241885, 1.318506389346667, 11.86655750412, 0.000911893273825466, 0.000911893273825466)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(5, $c).Value = $row5[$c - 1] }

$row6 = @("MuSCs", "Fgf18", "Fgfr1", "FAPs", 2.0, 0.6666666666666666, 0.71434, 2.14302, 0.06249742634325634, 0.06249742634325634, 3.0, 1.0, 82.95722966666666, 248.871689, 0.6557810310272387, 0.6557810310272387, 59.25966744008666, 533.33700696078, 0.04098462668392956, 0.04098462668392956)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(6, $c).Value = $row6[$c - 1] }

$row7 = @("MuSCs", "Fgf18", "Fgfr1", "MuSCs", 2.0, 0.6666666666666666, 0.71434, 2.14302, 0.06249742634325634, 0.06249742634325634, 3.0, 1.0, 41.69841866666667, 125.095256, 0.3296280757603424, 0.3296280757603424, 29.78684839034667, 268.08163551312, 0.02060090638550132, 0.02060090638550132)
for ($c = 1; $c -le 20; $c++) { $ws.Cells.Item(7, $c).Value = $row7[$c - 1] }
